$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column in H, matching the style already used by the
# other header cells (bold, centered, bordered -> same style index as B1).
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save indicator values for rows 2-7.
$saveValues = @(0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
